# Applies the crypto price/volume refresh for the "Tue Feb 21 09:28:31 UTC 2023"
# GitHub Actions update. Writes new Price (D) / Volume(1h) (E) values, and
# corrects the row ordering for a couple of coins whose rank changed
# (swap of Polygon/OKB in rows 9-10, and ShibaInu/Chainlink/WrappedEther
# in rows 15-17), by updating the Coin (B) and Link (C) text in place.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.880.09"
$ws.Range("E2").Value = "  +1.47%  "
$ws.Range("D3").Value = "1.687.99"
$ws.Range("E3").Value = "  -0.51%  "
$ws.Range("D4").Value = "1.006"
$ws.Range("E4").Value = "  +0.51%  "
$ws.Range("D5").Value = "315.54"
$ws.Range("E5").Value = "  -0.41%  "
$ws.Range("D6").Value = "1.006"
$ws.Range("E6").Value = "  +0.48%  "
$ws.Range("D7").Value = "0.3944"
$ws.Range("E7").Value = "  +0.62%  "
$ws.Range("D8").Value = "0.3975"
$ws.Range("E8").Value = "  -1.92%  "
$ws.Range("B9").Value = "OKB"
$ws.Range("C9").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D9").Value = "52.42"
$ws.Range("E9").Value = "  -2.81%  "
$ws.Range("B10").Value = "Polygon"
$ws.Range("C10").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D10").Value = "1.437"
$ws.Range("E10").Value = "  -3.23%  "
$ws.Range("D11").Value = "1.007"
$ws.Range("E11").Value = "  +0.45%  "
$ws.Range("D12").Value = "0.08718"
$ws.Range("E12").Value = "  -0.97%  "
$ws.Range("D13").Value = "25.46"
$ws.Range("E13").Value = "  -2.59%  "
$ws.Range("D14").Value = "7.370"
$ws.Range("E14").Value = "  -1.06%  "
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "2.084.58"
$ws.Range("E15").Value = "  +23.07%  "
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").Value = "0.00001335"
$ws.Range("E16").Value = "  -1.62%  "
$ws.Range("B17").Value = "Chainlink"
$ws.Range("C17").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D17").Value = "7.854"
$ws.Range("E17").Value = "  -3.12%  "
$ws.Range("D18").Value = "94.62"
$ws.Range("E18").Value = "  -3.06%  "
$ws.Range("D19").Value = "0.07246"
$ws.Range("E19").Value = "  +0.80%  "
$ws.Range("D20").Value = "20.40"
$ws.Range("E20").Value = "  -0.17%  "
$ws.Range("D21").Value = "7.161"
$ws.Range("E21").Value = "  -1.89%  "
$ws.Range("E22").Value = "  +0.59%  "
$ws.Range("D23").Value = "14.15"
$ws.Range("E23").Value = "  -1.28%  "
$ws.Range("D24").Value = "24.869.44"
$ws.Range("E24").Value = "  +1.47%  "
$ws.Range("D25").Value = "2.400"
$ws.Range("E25").Value = "  +2.89%  "
$ws.Range("D26").Value = "2.823"
$ws.Range("E26").Value = "  -6.40%  "
$ws.Range("D27").Value = "23.07"
$ws.Range("E27").Value = "  +0.50%  "
$ws.Range("D28").Value = "6.005"
$ws.Range("E28").Value = "  +1.65%  "
$ws.Range("D29").Value = "161.37"
$ws.Range("E29").Value = "  -4.49%  "
$ws.Range("D30").Value = "148.17"
$ws.Range("E30").Value = "  +2.24%  "
$ws.Range("D31").Value = "8.064"
$ws.Range("E31").Value = "  -4.00%  "
$ws.Range("D32").Value = "2.605"
$ws.Range("E32").Value = "  +19.40%  "
$ws.Range("D33").Value = "2.303.21"
$ws.Range("E33").Value = "  +22.52%  "
$ws.Range("D34").Value = "0.08502"
$ws.Range("E34").Value = "  -3.08%  "
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("D36").Value = "1.028"
$ws.Range("E36").Value = "  -1.89%  "
$ws.Range("D37").Value = "7.046"
$ws.Range("E37").Value = "  -2.12%  "
$ws.Range("D38").Value = "0.2855"
$ws.Range("E38").Value = "  +1.74%  "
$ws.Range("D39").Value = "0.09620"
$ws.Range("E39").Value = "  +4.83%  "
$ws.Range("D40").Value = "10.83"
$ws.Range("E40").Value = "  -0.47%  "
$ws.Range("D41").Value = "0.8068"
$ws.Range("E41").Value = "  -5.79%  "
$ws.Range("D42").Value = "13.90"
$ws.Range("E42").Value = "  -2.13%  "
$ws.Range("E43").Value = "  -0.85%  "
$ws.Range("D44").Value = "16.97"
$ws.Range("E44").Value = "  -2.57%  "
$ws.Range("D45").Value = "2.623"
$ws.Range("E45").Value = "  -1.45%  "
$ws.Range("D46").Value = "0.7261"
$ws.Range("E46").Value = "  -2.90%  "
$ws.Range("D47").Value = "4.218"
$ws.Range("E47").Value = "  -0.99%  "
$ws.Range("D48").Value = "0.08921"
$ws.Range("E48").Value = "  +8.52%  "
$ws.Range("D49").Value = "1.376"
$ws.Range("E49").Value = "  -0.91%  "
$ws.Range("E50").Value = "  +0.61%  "
$ws.Range("D51").Value = "139.84"
$ws.Range("E51").Value = "  -0.29%  "
